# The "mousse" row (row 10) is a blank placeholder row (only Column filled in,
# Description/Values empty, highlighted yellow) that needs to be removed
# entirely. Deleting the whole row shifts every row below it up by one,
# which is exactly what the target diff shows (rows 11-21 become rows 10-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(10).Delete()

# Move the selection to reflect where the cursor ended up after the edit.
$ws.Range("B18").Select()
